$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting cell value: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set the active cell / selection to E8
$ws.Range("E8").Select()
